$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$seleniumError = "no such element: Unable to locate element: {`"method`":`"xpath`",`"selector`":`"//button[@id='loginButton']`"}
  (Session info: chrome=134.0.6998.179)
For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Build info: version: '4.29.0', revision: '5fc1ec94cb'
System info: os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '23.0.1'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Command: [81ec64fa8daf30c077d70ca098674836, findElement {value=//button[@id='loginButton'], using=xpath}]
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 134.0.6998.179, chrome: {chromedriverVersion: 134.0.6998.165 (fd886e2cb29..., userDataDir: C:\Users\Admin\AppData\Loca...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:49661}, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:49661/devtoo..., se:cdpVersion: 134.0.6998.179, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 81ec64fa8daf30c077d70ca098674836"

# Insert a new row above row 13 ("Verify Home Page Loads Successfully") to make
# room for the new "resend otp" account-block test case, which shifts that last
# case down to row 14.
$ws.Rows.Item(13).Insert()

# Row 2: first test case stays PASSED, only its execution time moved to the new run.
$ws.Cells.Item(2, 4).Value = "03/04/2025 12:07:58 PM"

# Row 3: "Verify empty email state" now FAILED (couldn't find the login button),
# with the full Selenium/ChromeDriver error captured as the comment.
$ws.Cells.Item(3, 3).Value = "FAILED"
$ws.Cells.Item(3, 4).Value = "03/04/2025 12:08:03 PM"
$ws.Cells.Item(3, 5).Value = $seleniumError

# Rows 4-8: unchanged test case names, but all SKIPPED as a consequence of the
# earlier failure, each stamped with the later run's timestamp.
for ($r = 4; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = "SKIPPED"
    $ws.Cells.Item($r, 4).Value = "03/04/2025 12:08:03 PM"
    $ws.Cells.Item($r, 5).Value = "Test was skipped."
}

# Rows 9-11: these three test cases moved up one row (their old slots were at
# rows 10-12), keeping the same relative order, all SKIPPED.
$ws.Cells.Item(9, 2).Value = "Verify account block after attempting wrong OTP for 5 times"
$ws.Cells.Item(10, 2).Value = "Verify Go To Sign In page Navigation"
$ws.Cells.Item(11, 2).Value = "Verify that navigation and getOTP blocked for blocked account"
for ($r = 9; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = "SKIPPED"
    $ws.Cells.Item($r, 4).Value = "03/04/2025 12:08:03 PM"
    $ws.Cells.Item($r, 5).Value = "Test was skipped."
}

# Row 12: "Verify Resend OTP button" moved down to directly precede the new
# resend-otp block-account test case; also SKIPPED in this run.
$ws.Cells.Item(12, 2).Value = "Verify Resend OTP button"
$ws.Cells.Item(12, 3).Value = "SKIPPED"
$ws.Cells.Item(12, 4).Value = "03/04/2025 12:08:03 PM"
$ws.Cells.Item(12, 5).Value = "Test was skipped."

# Row 13 (newly inserted): brand-new test case covering the resend-otp block issue.
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Verify account block after 5 times of resend otp"
$ws.Cells.Item(13, 3).Value = "SKIPPED"
$ws.Cells.Item(13, 4).Value = "03/04/2025 12:08:03 PM"
$ws.Cells.Item(13, 5).Value = "Test was skipped."

# Row 14 (shifted down from the old row 13): "Home Page Loads Successfully", PASSED again.
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Verify Home Page Loads Successfully"
$ws.Cells.Item(14, 3).Value = "PASSED"
$ws.Cells.Item(14, 4).Value = "03/04/2025 12:08:03 PM"
$ws.Cells.Item(14, 5).Value = "Test executed successfully."
